# Update the ASPP variable crosswalk: add a new "var_2017" column (G) that
# mirrors the existing "var_2016" column (F) for rows 1-36, matches its
# column formatting, and leaves the worksheet with frozen header/first
# column panes and the newly-entered range selected (as an analyst would
# leave it after pasting in the new year's column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the added column.
$ws.Cells.Item(1, 7).Value = "var_2017"

# Mirror column F ("var_2016") into the new column G ("var_2017") for the
# rows that already have data (rows 2-36).
for ($r = 2; $r -le 36; $r++) {
    $srcValue = $ws.Cells.Item($r, 6).Value()
    $ws.Cells.Item($r, 7).Value = $srcValue
}

# Match the column width/formatting that columns B:F already share.
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(6).ColumnWidth()

# Freeze the header row and the first (vname) column, then leave the
# freshly-filled cell selected in the lower-right pane.
[void]$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
[void]$ws.Range("G5").Select()
